$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 31.27132566666667
$ws.Range("H2").Value = 93.81397700000001
$ws.Range("I2").Value = 0.9493361071405608
$ws.Range("J2").Value = 0.9493361071405608
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.11318066666666
$ws.Range("N2").Value = 66.33954199999999
$ws.Range("O2").Value = 0.1993293533530854
$ws.Range("P2").Value = 0.1993293533530854
$ws.Range("Q2").Value = 691.5084741531705
$ws.Range("R2").Value = 6223.576267378534
$ws.Range("S2").Value = 0.1892305523510634
$ws.Range("T2").Value = 0.1892305523510634

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 31.27132566666667
$ws.Range("H3").Value = 93.81397700000001
$ws.Range("I3").Value = 0.9493361071405608
$ws.Range("J3").Value = 0.9493361071405608
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.622575
$ws.Range("N3").Value = 22.867725
$ws.Range("O3").Value = 0.06871028498970018
$ws.Range("P3").Value = 0.06871028498970018
$ws.Range("Q3").Value = 238.3680252435917
$ws.Range("R3").Value = 2145.312227192325
$ws.Range("S3").Value = 0.06522915447264048
$ws.Range("T3").Value = 0.06522915447264048

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 31.27132566666667
$ws.Range("H4").Value = 93.81397700000001
$ws.Range("I4").Value = 0.9493361071405608
$ws.Range("J4").Value = 0.9493361071405608
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 81.20214833333334
$ws.Range("N4").Value = 243.606445
$ws.Range("O4").Value = 0.7319603616572145
$ws.Range("P4").Value = 0.7319603616572145
$ws.Range("Q4").Value = 2539.298825364641
$ws.Range("R4").Value = 22853.68942828177
$ws.Range("S4").Value = 0.6948764003168569
$ws.Range("T4").Value = 0.6948764003168569

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.96805
$ws.Range("H5").Value = 2.90415
$ws.Range("I5").Value = 0.02938809912676722
$ws.Range("J5").Value = 0.02938809912676721
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.11318066666666
$ws.Range("N5").Value = 66.33954199999999
$ws.Range("O5").Value = 0.1993293533530854
$ws.Range("P5").Value = 0.1993293533530854
$ws.Range("Q5").Value = 21.40666454436666
$ws.Range("R5").Value = 192.6599808993
$ws.Range("S5").Value = 0.005857910795214883
$ws.Range("T5").Value = 0.005857910795214882

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.96805
$ws.Range("H6").Value = 2.90415
$ws.Range("I6").Value = 0.02938809912676722
$ws.Range("J6").Value = 0.02938809912676721
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.622575
$ws.Range("N6").Value = 22.867725
$ws.Range("O6").Value = 0.06871028498970018
$ws.Range("P6").Value = 0.06871028498970018
$ws.Range("Q6").Value = 7.37903372875
$ws.Range("R6").Value = 66.41130355875001
$ws.Range("S6").Value = 0.002019264666305734
$ws.Range("T6").Value = 0.002019264666305734

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.96805
$ws.Range("H7").Value = 2.90415
$ws.Range("I7").Value = 0.02938809912676722
$ws.Range("J7").Value = 0.02938809912676721
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 81.20214833333334
$ws.Range("N7").Value = 243.606445
$ws.Range("O7").Value = 0.7319603616572145
$ws.Range("P7").Value = 0.7319603616572145
$ws.Range("Q7").Value = 78.60773969408334
$ws.Range("R7").Value = 707.46965724675
$ws.Range("S7").Value = 0.0215109236652466
$ws.Range("T7").Value = 0.0215109236652466

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.700829
$ws.Range("H8").Value = 2.102487
$ws.Range("I8").Value = 0.02127579373267201
$ws.Range("J8").Value = 0.02127579373267201
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 22.11318066666666
$ws.Range("N8").Value = 66.33954199999999
$ws.Range("O8").Value = 0.1993293533530854
$ws.Range("P8").Value = 0.1993293533530854
$ws.Range("Q8").Value = 15.49755829343933
$ws.Range("R8").Value = 139.478024640954
$ws.Range("S8").Value = 0.004240890206807139
$ws.Range("T8").Value = 0.004240890206807139

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.700829
$ws.Range("H9").Value = 2.102487
$ws.Range("I9").Value = 0.02127579373267201
$ws.Range("J9").Value = 0.02127579373267201
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.622575
$ws.Range("N9").Value = 22.867725
$ws.Range("O9").Value = 0.06871028498970018
$ws.Range("P9").Value = 0.06871028498970018
$ws.Range("Q9").Value = 5.342121614675
$ws.Range("R9").Value = 48.079094532075
$ws.Range("S9").Value = 0.001461865850753971
$ws.Range("T9").Value = 0.001461865850753971

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.700829
$ws.Range("H10").Value = 2.102487
$ws.Range("I10").Value = 0.02127579373267201
$ws.Range("J10").Value = 0.02127579373267201
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 81.20214833333334
$ws.Range("N10").Value = 243.606445
$ws.Range("O10").Value = 0.7319603616572145
$ws.Range("P10").Value = 0.7319603616572145
$ws.Range("Q10").Value = 56.90882041430167
$ws.Range("R10").Value = 512.179383728715
$ws.Range("S10").Value = 0.0155730376751109
$ws.Range("T10").Value = 0.0155730376751109
